$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# 1) Update the "Valor Mora" total (E11) to reflect the newly added period.
# 2) Update "Cant. Periodos" (F13) from 1 to 2 (now covering periods 2507 and
#    2508).
# 3) Insert a brand-new block of 3 worker rows for period 2508 (a duplicate
#    of the 3 existing workers for period 2507), placed right after the
#    existing table, pushing the signature/footer block further down.
# ---------------------------------------------------------------------------

# --- Totals header -----------------------------------------------------
$ws.Range("E11").Value = 367760
$ws.Range("F13").Value = 2

# --- Make room for 3 new worker rows right after the current last data
#     row (row 18), before the blank spacer + footer rows.
$ws.Rows("19:21").Insert()

# New rows 19/20 should look exactly like the two "interior" worker rows
# (16/17); row 21 should look like the old "closing" worker row (18, which
# carries the thicker bottom border that closes the table) since it is now
# the new last row of the table.
$ws.Range("B16:J16").Copy($ws.Range("B19:J19"))
$ws.Range("B17:J17").Copy($ws.Range("B20:J20"))
$ws.Range("B18:J18").Copy($ws.Range("B21:J21"))

# The old row 18 (previously the closing row of the table) is no longer the
# last row, so it now gets the plain "interior" row formatting instead (same
# as rows 16/17).
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# --- Row 18: SANDRA MILENA ALDANA RODRIGUEZ, period 2507 (unchanged data,
#     now with interior-row styling) ------------------------------------
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047417283"
$ws.Range("D18").Value = "SANDRA MILENA ALDANA RODRIGUEZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 70000
$ws.Range("G18").Value = 1750000

# --- Row 19: LUIS GABRIEL WATTS PAJARO, period 2508 ---------------------
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73576599"
$ws.Range("D19").Value = "LUIS GABRIEL WATTS PAJARO"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- Row 20: MARIA CAROLINA LORDUY ALCALA, period 2508 ------------------
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45556506"
$ws.Range("D20").Value = "MARIA CAROLINA LORDUY ALCALA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# --- Row 21: SANDRA MILENA ALDANA RODRIGUEZ, period 2508 (new closing
#     row of the table) ---------------------------------------------------
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047417283"
$ws.Range("D21").Value = "SANDRA MILENA ALDANA RODRIGUEZ"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 70000
$ws.Range("G21").Value = 1750000
